$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff/Handback datetimes
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 04:21:50"
$wsZhCn.Range("E3").Value = "2016-03-22 04:21:50"
$wsZhCn.Range("H2").Value = "2016-03-22 04:22:13"
$wsZhCn.Range("H3").Value = "2016-03-22 04:22:13"

# de-de sheet: update Correspond Handoff/Handback datetimes
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 04:21:53"
$wsDeDe.Range("E3").Value = "2016-03-22 04:21:53"
$wsDeDe.Range("H2").Value = "2016-03-22 04:22:19"
$wsDeDe.Range("H3").Value = "2016-03-22 04:22:19"
